# save valor con su funcion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - function names / parameter labels
$ws.Range("A5").Value = "validarnumero"
$ws.Range("B5").Value = "valoresposibles"
$ws.Range("C5").Value = "validarnumero"
$ws.Range("D5").Value = "validarcaracteres"
$ws.Range("E5").Value = "valoresposilbes"
$ws.Range("F5").Value = "validarnumero"

# Row 6 - extra parameter
$ws.Range("C6").Value = "longitudadmaxima"

# Row 8 - value/type descriptions
$ws.Range("A8").Value = "Valor – string"
$ws.Range("B8").Value = "valor string, valoresPosibles [] string"
$ws.Range("C8").Value = "Valor – string, "
$ws.Range("D8").Value = "texto string"
$ws.Range("E8").Value = "valor string, valoresPosibles [] string"
$ws.Range("F8").Value = "Valor – string"

# Row 9 - extra description
$ws.Range("C9").Value = "Valor – string, longitud int"

# Column widths (values chosen so the saved OOXML <col width=.../> lands on
# 13.38 / 31.79 / 22.97 / 15.47 / 31.79 / 13.38 as closely as the engine's
# pixel-snapping allows)
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(2).ColumnWidth = 31.0
$ws.Columns.Item(3).ColumnWidth = 22.166666666666668
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666
$ws.Columns.Item(5).ColumnWidth = 31.0
$ws.Columns.Item(6).ColumnWidth = 12.5

# Selection
$ws.Range("B2").Select()
